$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column I (sum_fiber). Columns J (sum_fats) and K (sum_sodium)
# shift left to become I and J respectively.
$ws.Columns.Item(9).Delete()
